$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F17").Value = 1774
$ws.Range("F18").Value = 3494
$ws.Range("F24").Value = 3362
$ws.Range("F35").Value = 1208
$ws.Range("F36").Value = 1869
$ws.Range("F40").Value = 237
$ws.Range("F43").Value = 65

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 16

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 203

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F16").Value = 1774
$ws.Range("F17").Value = 3494
$ws.Range("F22").Value = 3362
$ws.Range("F28").Value = 16
$ws.Range("F39").Value = 1208
$ws.Range("F40").Value = 1869
$ws.Range("F45").Value = 237
$ws.Range("F48").Value = 65
